$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 21 (shifts rows 21-28 down to 22-29,
#    and shifts the picture anchored below it down as well).
$ws.Rows(21).Insert()

# 2. Populate the newly inserted row's O21 cell with the new note text.
$ws.Range("O21").Value = "Nếu soi trên khung H4 chắc chắn phải có ít nhất một cây nến vàng được sinh ra, nếu không có hãy cẩn thận"

# 3. Move the picture down by the height of the newly inserted row so it
#    keeps sitting on the same rows relative to the text above it. The
#    picture used to start at (1-indexed) row 30 with a 6024 EMU offset
#    into that row; after the insert it belongs in row 31, same offset.
#    We recompute the target using the row's Top (in points) plus the
#    original sub-row EMU offset converted to points, rather than
#    reading back Shape.Top (which the host only reports rounded to
#    2 decimals and would otherwise introduce drift).
$pic = $ws.Shapes.Item(1)
$rowOffPt = 6024 / 12700.0
$pic.Top = $ws.Rows(31).Top + $rowOffPt

# 4. Update the view: scroll back to the top (removes the stale
#    topLeftCell="A43") and move the active selection to J24.
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J24").Select()
